$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.282228905847376
$ws.Range("C2").Value2 = 0.5481659376659422
$ws.Range("D2").Value2 = 0.04427475039713036
$ws.Range("E2").Value2 = 0.09086768655857469
$ws.Range("F2").Value2 = 3.156006502200938
$ws.Range("H2").Value2 = 0.07973214163530429
$ws.Range("I2").Value2 = 1.567393841450439
$ws.Range("J2").Value2 = 0.1622475285560512
$ws.Range("M2").Value2 = 0.5066399124069605
$ws.Range("B3").Value2 = 1.200853625392824
$ws.Range("C3").Value2 = 0.513979175820225
$ws.Range("D3").Value2 = 0.04430237240105983
$ws.Range("E3").Value2 = 0.09104512587845726
$ws.Range("F3").Value2 = 3.109600644876238
$ws.Range("H3").Value2 = 0.07973214163530429
$ws.Range("I3").Value2 = 1.549139059382057
$ws.Range("J3").Value2 = 0.161989779380832
$ws.Range("M3").Value2 = 0.4871025092023942
$ws.Range("B4").Value2 = 1.151708354095604
$ws.Range("C4").Value2 = 0.4933472951890963
$ws.Range("D4").Value2 = 0.04432452183462798
$ws.Range("E4").Value2 = 0.0911834738309274
$ws.Range("F4").Value2 = 3.082774053227809
$ws.Range("H4").Value2 = 0.07973214163530429
$ws.Range("I4").Value2 = 1.538704150164619
$ws.Range("J4").Value2 = 0.1618953660720734
$ws.Range("M4").Value2 = 0.4754341351725273
$ws.Range("B5").Value2 = 1.131886617910197
$ws.Range("C5").Value2 = 0.4850291901940977
$ws.Range("D5").Value2 = 0.04433485433509254
$ws.Range("E5").Value2 = 0.09124725632208452
$ws.Range("F5").Value2 = 3.072258635838722
$ws.Range("H5").Value2 = 0.07973214163530429
$ws.Range("I5").Value2 = 1.534645015360226
$ws.Range("J5").Value2 = 0.1618729085122794
$ws.Range("M5").Value2 = 0.4707613769008105
$ws.Range("B6").Value2 = 1.128607617262958
$ws.Range("C6").Value2 = 0.4836533638362255
$ws.Range("D6").Value2 = 0.04433664899382039
$ws.Range("E6").Value2 = 0.09125829492779758
$ws.Range("F6").Value2 = 3.070537650114403
$ws.Range("H6").Value2 = 0.07973214163530429
$ws.Range("I6").Value2 = 1.533982626986955
$ws.Range("J6").Value2 = 0.1618701457202008
$ws.Range("M6").Value2 = 0.4699904280438219
$ws.Range("B7").Value2 = 1.151440200643378
$ws.Range("C7").Value2 = 0.4932347526287515
$ws.Range("D7").Value2 = 0.04432465588946011
$ws.Range("E7").Value2 = 0.09118430402687849
$ws.Range("F7").Value2 = 3.082630555006318
$ws.Range("H7").Value2 = 0.07973214163530429
$ws.Range("I7").Value2 = 1.538648626900354
$ws.Range("J7").Value2 = 0.1618949983999158
$ws.Range("M7").Value2 = 0.4753707841543502
$ws.Range("B8").Value2 = 1.254000162951058
$ws.Range("C8").Value2 = 0.5363034081224214
$ws.Range("D8").Value2 = 0.04428319819928106
$ws.Range("E8").Value2 = 0.09092277268185001
$ws.Range("F8").Value2 = 3.139658294903171
$ws.Range("H8").Value2 = 0.07973214163530429
$ws.Range("I8").Value2 = 1.56093820798759
$ws.Range("J8").Value2 = 0.1621453821328984
$ws.Range("M8").Value2 = 0.499835278947856
$ws.Range("B9").Value2 = 1.461667651512073
$ws.Range("C9").Value2 = 0.6236460041066607
$ws.Range("D9").Value2 = 0.04424303362943949
$ws.Range("E9").Value2 = 0.09064272450688726
$ws.Range("F9").Value2 = 3.264836396595712
$ws.Range("H9").Value2 = 0.07973214163530429
$ws.Range("I9").Value2 = 1.610850791761834
$ws.Range("J9").Value2 = 0.1631448726830129
$ws.Range("M9").Value2 = 0.5504225436411758
$ws.Range("B10").Value2 = 1.618315520527176
$ws.Range("C10").Value2 = 0.6896348235306959
$ws.Range("D10").Value2 = 0.04423858127548308
$ws.Range("E10").Value2 = 0.09057838001962715
$ws.Range("F10").Value2 = 3.365123124212118
$ws.Range("H10").Value2 = 0.07973214163530429
$ws.Range("I10").Value2 = 1.651397602501632
$ws.Range("J10").Value2 = 0.164192099401653
$ws.Range("M10").Value2 = 0.5892043350370599
$ws.Range("B11").Value2 = 1.690483181500724
$ws.Range("C11").Value2 = 0.7200632040287474
$ws.Range("D11").Value2 = 0.04424200005815093
$ws.Range("E11").Value2 = 0.09057971415335153
$ws.Range("F11").Value2 = 3.412592780672639
$ws.Range("H11").Value2 = 0.07973214163530429
$ws.Range("I11").Value2 = 1.670705889722868
$ws.Range("J11").Value2 = 0.1647371099159898
$ws.Range("M11").Value2 = 0.6072033333615821
$ws.Range("B12").Value2 = 1.717943016845311
$ws.Range("C12").Value2 = 0.7316455331151701
$ws.Range("D12").Value2 = 0.04424407802852315
$ws.Range("E12").Value2 = 0.09058461136797824
$ws.Range("F12").Value2 = 3.430837209423061
$ws.Range("H12").Value2 = 0.07973214163530429
$ws.Range("I12").Value2 = 1.678143188681943
$ws.Range("J12").Value2 = 0.1649534087570501
$ws.Range("M12").Value2 = 0.6140707741688658
$ws.Range("B13").Value2 = 1.712023185112059
$ws.Range("C13").Value2 = 0.729148399760561
$ws.Range("D13").Value2 = 0.04424359564932967
$ws.Range("E13").Value2 = 0.09058336146783752
$ws.Range("F13").Value2 = 3.426895948897482
$ws.Range("H13").Value2 = 0.07973214163530429
$ws.Range("I13").Value2 = 1.676535819698373
$ws.Range("J13").Value2 = 0.1649063831865902
$ws.Range("M13").Value2 = 0.6125894479265952
$ws.Range("B14").Value2 = 1.692739677314989
$ws.Range("C14").Value2 = 0.7210148865996757
$ws.Range("D14").Value2 = 0.04424215531068065
$ws.Range("E14").Value2 = 0.09058002907770657
$ws.Range("F14").Value2 = 3.414088359177583
$ws.Range("H14").Value2 = 0.07973214163530429
$ws.Range("I14").Value2 = 1.671315233410567
$ws.Range("J14").Value2 = 0.1647547059623804
$ws.Range("M14").Value2 = 0.6077672858467906
$ws.Range("B15").Value2 = 1.680945127950906
$ws.Range("C15").Value2 = 0.7160406830375905
$ws.Range("D15").Value2 = 0.04424137509662884
$ws.Range("E15").Value2 = 0.09057855960245043
$ws.Range("F15").Value2 = 3.406278418684877
$ws.Range("H15").Value2 = 0.07973214163530429
$ws.Range("I15").Value2 = 1.66813388487374
$ws.Range("J15").Value2 = 0.1646630919546936
$ws.Range("M15").Value2 = 0.6048203020945238
$ws.Range("B16").Value2 = 1.613617542779309
$ws.Range("C16").Value2 = 0.687654573766622
$ws.Range("D16").Value2 = 0.04423846743501869
$ws.Range("E16").Value2 = 0.09057890790669454
$ws.Range("F16").Value2 = 3.36205831491327
$ws.Range("H16").Value2 = 0.07973214163530429
$ws.Range("I16").Value2 = 1.650153262963173
$ws.Range("J16").Value2 = 0.1641578665666401
$ws.Range("M16").Value2 = 0.5880352648396894
$ws.Range("B17").Value2 = 1.572547427185327
$ws.Range("C17").Value2 = 0.6703461511021374
$ws.Range("D17").Value2 = 0.04423807830435411
$ws.Range("E17").Value2 = 0.09058695373831505
$ws.Range("F17").Value2 = 3.335406237059289
$ws.Range("H17").Value2 = 0.07973214163530429
$ws.Range("I17").Value2 = 1.639344934870451
$ws.Range("J17").Value2 = 0.1638655363569015
$ws.Range("M17").Value2 = 0.5778298106269375
$ws.Range("B18").Value2 = 1.549010364526339
$ws.Range("C18").Value2 = 0.6604293188091219
$ws.Range("D18").Value2 = 0.04423836686897165
$ws.Range("E18").Value2 = 0.0905944626656936
$ws.Range("F18").Value2 = 3.320250504437325
$ws.Range("H18").Value2 = 0.07973214163530429
$ws.Range("I18").Value2 = 1.633209440244968
$ws.Range("J18").Value2 = 0.1637038518788003
$ws.Range("M18").Value2 = 0.5719934852050628
$ws.Range("B19").Value2 = 1.541055754982892
$ws.Range("C19").Value2 = 0.65707823938709
$ws.Range("D19").Value2 = 0.04423855256579623
$ws.Range("E19").Value2 = 0.09059750018187707
$ws.Range("F19").Value2 = 3.315148804056179
$ws.Range("H19").Value2 = 0.07973214163530429
$ws.Range("I19").Value2 = 1.631145961996253
$ws.Range("J19").Value2 = 0.1636502157054593
$ws.Range("M19").Value2 = 0.5700231651225494
$ws.Range("B20").Value2 = 1.576910567530831
$ws.Range("C20").Value2 = 0.6721846700891092
$ws.Range("D20").Value2 = 0.04423806669504216
$ws.Range("E20").Value2 = 0.09058579912059095
$ws.Range("F20").Value2 = 3.338225384131647
$ws.Range("H20").Value2 = 0.07973214163530429
$ws.Range("I20").Value2 = 1.640487088095298
$ws.Range("J20").Value2 = 0.1638959868816414
$ws.Range("M20").Value2 = 0.5789127210671126
$ws.Range("B21").Value2 = 1.698400133904954
$ws.Range("C21").Value2 = 0.7234022678422889
$ws.Range("D21").Value2 = 0.04424255710822744
$ws.Range("E21").Value2 = 0.09058088875302595
$ws.Range("F21").Value2 = 3.417842943363553
$ws.Range("H21").Value2 = 0.07973214163530429
$ws.Range("I21").Value2 = 1.672845223313445
$ws.Range("J21").Value2 = 0.1647989877617348
$ws.Range("M21").Value2 = 0.6091822687481709
$ws.Range("B22").Value2 = 1.778567605100022
$ws.Range("C22").Value2 = 0.7572246718000883
$ws.Range("D22").Value2 = 0.0442500582821852
$ws.Range("E22").Value2 = 0.09060327562829684
$ws.Range("F22").Value2 = 3.471445160357774
$ws.Range("H22").Value2 = 0.07973214163530429
$ws.Range("I22").Value2 = 1.694726350248871
$ws.Range("J22").Value2 = 0.1654469669525866
$ws.Range("M22").Value2 = 0.6292661339806784
$ws.Range("B23").Value2 = 1.735710201098982
$ws.Range("C23").Value2 = 0.7391408245700859
$ws.Range("D23").Value2 = 0.04424563669012738
$ws.Range("E23").Value2 = 0.09058898814592098
$ws.Range("F23").Value2 = 3.442692265771655
$ws.Range("H23").Value2 = 0.07973214163530429
$ws.Range("I23").Value2 = 1.682980385206648
$ws.Range("J23").Value2 = 0.1650958222088548
$ws.Range("M23").Value2 = 0.6185193630266497
$ws.Range("B24").Value2 = 1.574937759198292
$ws.Range("C24").Value2 = 0.671353369782878
$ws.Range("D24").Value2 = 0.04423807034798299
$ws.Range("E24").Value2 = 0.09058631214195323
$ws.Range("F24").Value2 = 3.336950328093394
$ws.Range("H24").Value2 = 0.07973214163530429
$ws.Range("I24").Value2 = 1.639970476672161
$ws.Range("J24").Value2 = 0.1638822003305336
$ws.Range("M24").Value2 = 0.5784230409499003
$ws.Range("B25").Value2 = 1.404778708955803
$ws.Range("C25").Value2 = 0.5997025805874614
$ws.Range("D25").Value2 = 0.04424949877610285
$ws.Range("E25").Value2 = 0.09069361587740765
$ws.Range("F25").Value2 = 3.229524848490968
$ws.Range("H25").Value2 = 0.07973214163530429
$ws.Range("I25").Value2 = 1.596674322333826
$ws.Range("J25").Value2 = 0.1628197483950728
$ws.Range("M25").Value2 = 0.5364553647189467
